$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = 0
$ws.Range("C1").Value = "Monday"
$ws.Range("D1").Value = "Tuesday"
$ws.Range("E1").Value = "Wednesday"
$ws.Range("F1").Value = "Thursday"
$ws.Range("G1").Value = "Friday"
$ws.Range("H1").Value = "Saturday"
$ws.Range("I1").Value = "Sunday"
$ws.Range("B1").Value = "WholeWeek"
$ws.Range("J1").WrapText = $true
